$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2453
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 60
$ws.Range("G2").Value = 14
$ws.Range("H2").Value = 18
$ws.Range("I2").Value = 18
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2143
$ws.Range("L2").Value = 1237
$ws.Range("M2").Value = 905
$ws.Range("N2").Value = 905
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 290
$ws.Range("Q2").Value = 206
$ws.Range("R2").Value = -102
$ws.Range("S2").Value = -119
$ws.Range("T2").Value = 87
$ws.Range("U2").Value = 118
$ws.Range("V2").Value = 867
$ws.Range("W2").Value = 2.44
$ws.Range("X2").Value = 0.75
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 0.85
$ws.Range("AA2").Value = 136.66
$ws.Range("AB2").Value = 211.89
$ws.Range("AC2").Value = 59
$ws.Range("AD2").Value = 77.58
$ws.Range("AE2").Value = 2981
$ws.Range("AF2").Value = 1.54
$ws.Range("AG2").Value = 48
$ws.Range("AH2").Value = 1.04
$ws.Range("AI2").Value = 80.90000000000001
$ws.Range("AJ2").Value = 30346104

# Row 3
$ws.Range("D3").Value = 2504
$ws.Range("E3").Value = 142
$ws.Range("F3").Value = 142
$ws.Range("G3").Value = 105
$ws.Range("H3").Value = 90
$ws.Range("I3").Value = 89
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2086
$ws.Range("L3").Value = 1111
$ws.Range("M3").Value = 975
$ws.Range("N3").Value = 974
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 290
$ws.Range("Q3").Value = 242
$ws.Range("R3").Value = -64
$ws.Range("S3").Value = -170
$ws.Range("T3").Value = 54
$ws.Range("U3").Value = 188
$ws.Range("V3").Value = 745
$ws.Range("W3").Value = 5.67
$ws.Range("X3").Value = 3.59
$ws.Range("Y3").Value = 9.529999999999999
$ws.Range("Z3").Value = 4.25
$ws.Range("AA3").Value = 113.93
$ws.Range("AB3").Value = 235.76
$ws.Range("AC3").Value = 295
$ws.Range("AD3").Value = 23.49
$ws.Range("AE3").Value = 3209
$ws.Range("AF3").Value = 2.16
$ws.Range("AG3").Value = 72
$ws.Range("AH3").Value = 1.03
$ws.Range("AI3").Value = 24.31
$ws.Range("AJ3").Value = 30346104

# Row 4
$ws.Range("D4").Value = 2384
$ws.Range("E4").Value = 179
$ws.Range("F4").Value = 179
$ws.Range("G4").Value = 215
$ws.Range("H4").Value = 177
$ws.Range("I4").Value = 176
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3083
$ws.Range("L4").Value = 1335
$ws.Range("M4").Value = 1748
$ws.Range("N4").Value = 1747
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 380
$ws.Range("Q4").Value = 382
$ws.Range("R4").Value = -712
$ws.Range("S4").Value = 439
$ws.Range("T4").Value = 263
$ws.Range("U4").Value = 119
$ws.Range("V4").Value = 577
$ws.Range("W4").Value = 7.52
$ws.Range("X4").Value = 7.42
$ws.Range("Y4").Value = 12.97
$ws.Range("Z4").Value = 6.84
$ws.Range("AA4").Value = 76.34999999999999
$ws.Range("AB4").Value = 367.58
$ws.Range("AC4").Value = 514
$ws.Range("AD4").Value = 12.76
$ws.Range("AE4").Value = 4596
$ws.Range("AF4").Value = 1.43
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 1.52
$ws.Range("AI4").Value = 21.54
$ws.Range("AJ4").Value = 38000000

# Row 5
$ws.Range("D5").Value = 2622
$ws.Range("E5").Value = 196
$ws.Range("F5").Value = 196
$ws.Range("G5").Value = 236
$ws.Range("H5").Value = 175
$ws.Range("I5").Value = 175
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3538
$ws.Range("L5").Value = 1593
$ws.Range("M5").Value = 1945
$ws.Range("N5").Value = 1944
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 380
$ws.Range("Q5").Value = 429
$ws.Range("R5").Value = -583
$ws.Range("S5").Value = 136
$ws.Range("T5").Value = 480
$ws.Range("U5").Value = -51
$ws.Range("V5").Value = 766
$ws.Range("W5").Value = 7.48
$ws.Range("X5").Value = 6.68
$ws.Range("Y5").Value = 9.48
$ws.Range("Z5").Value = 5.29
$ws.Range("AA5").Value = 81.89
$ws.Range("AB5").Value = 402.38
$ws.Range("AC5").Value = 460
$ws.Range("AD5").Value = 17.03
$ws.Range("AE5").Value = 5114
$ws.Range("AF5").Value = 1.53
$ws.Range("AG5").Value = 120
$ws.Range("AH5").Value = 1.53
$ws.Range("AI5").Value = 26.06
$ws.Range("AJ5").Value = 38000000

# Row 6
$ws.Range("D6").Value = 3213
$ws.Range("E6").Value = 194
$ws.Range("F6").Value = 194
$ws.Range("G6").Value = 162
$ws.Range("H6").Value = 123
$ws.Range("I6").Value = 123
$ws.Range("J6").ClearContents()
$ws.Range("K6").Value = 3410
$ws.Range("L6").Value = 1436
$ws.Range("M6").Value = 1973
$ws.Range("N6").Value = 1972
$ws.Range("O6").ClearContents()
$ws.Range("P6").Value = 380
$ws.Range("Q6").Value = 343
$ws.Range("R6").Value = -181
$ws.Range("S6").Value = -185
$ws.Range("T6").Value = 152
$ws.Range("U6").Value = 192
$ws.Range("V6").Value = 642
$ws.Range("W6").Value = 6.03
$ws.Range("X6").Value = 3.82
$ws.Range("Y6").Value = 6.28
$ws.Range("Z6").Value = 3.53
$ws.Range("AA6").Value = 72.78
$ws.Range("AB6").Value = 419.26
$ws.Range("AC6").Value = 324
$ws.Range("AD6").Value = 17.34
$ws.Range("AE6").Value = 5189
$ws.Range("AF6").Value = 1.08
$ws.Range("AG6").Value = 120
$ws.Range("AH6").Value = 2.14
$ws.Range("AI6").Value = 37.09
$ws.Range("AJ6").Value = 38000000

# Rows 7-9: clear all data cells, keep A/B/C
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()